$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 22974.75
$ws.Range("J7").Value = 22974.75
$ws.Range("L7").Value = 22974.75
$ws.Range("N7").Value = -23198.75
$ws.Range("H14").Value = 22974.75
$ws.Range("J14").Value = 22974.75
$ws.Range("L14").Value = 22974.75
$ws.Range("N14").Value = -23356.75
$ws.Range("H21").Value = 25300
$ws.Range("I21").Value = 2950
$ws.Range("K21").Value = 2950
$ws.Range("M21").Value = -2482
$ws.Range("H23").Value = 25300
$ws.Range("I23").Value = 2950
$ws.Range("K23").Value = 2950
$ws.Range("M23").Value = -2716
$ws.Range("H44").Value = 50000
$ws.Range("J44").Value = 50000
$ws.Range("L44").Value = 50000
$ws.Range("N44").Value = -50924
$ws.Range("H74").Value = 10005879
$ws.Range("I74").Value = 16670665
$ws.Range("J74").Value = 8701
$ws.Range("K74").Value = 16670665
$ws.Range("L74").Value = 8701
$ws.Range("M74").Value = -16669729
$ws.Range("N74").Value = -10573
$ws.Range("H77").Value = 10005879
$ws.Range("I77").Value = 16670665
$ws.Range("J77").Value = 8701
$ws.Range("K77").Value = 83353325
$ws.Range("L77").Value = 43505
$ws.Range("M77").Value = -83348645
$ws.Range("N77").Value = -52865
$ws.Range("H129").Value = 870.96844
$ws.Range("I129").Value = 403.75
$ws.Range("J129").Value = 891.5055
$ws.Range("K129").Value = 1211.25
$ws.Range("L129").Value = 2674.5165
$ws.Range("M129").Value = 3788.75
$ws.Range("N129").Value = -12674.5165
$ws.Range("H138").Value = 2397.71
$ws.Range("I138").Value = 662.4
$ws.Range("J138").Value = 2831.5376
$ws.Range("K138").Value = 1987.2
$ws.Range("L138").Value = 8494.612800000001
$ws.Range("M138").Value = 3152.8
$ws.Range("N138").Value = -18774.6128
$ws.Range("H141").Value = 168466.08
$ws.Range("I141").Value = 223465.89
$ws.Range("K141").Value = 670397.67
$ws.Range("M141").Value = -665217.67

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3759.1428
$ws.Range("I122").Value = 1350
$ws.Range("K122").Value = 4050
$ws.Range("M122").Value = -1600
$ws.Range("H139").Value = 42781.07
$ws.Range("J139").Value = 42781.07
$ws.Range("L139").Value = 42781.07
$ws.Range("N139").Value = -53061.07

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 975.8
$ws.Range("I94").Value = 745.1667
$ws.Range("J94").Value = 1898.3334
$ws.Range("K94").Value = 745.1667
$ws.Range("L94").Value = 1898.3334
$ws.Range("M94").Value = -294.1667
$ws.Range("N94").Value = -2800.3334
$ws.Range("H134").Value = 2578.28
$ws.Range("I134").Value = 1710.4667
$ws.Range("J134").Value = 3880
$ws.Range("K134").Value = 5131.4001
$ws.Range("L134").Value = 11640
$ws.Range("M134").Value = -2596.4001
$ws.Range("N134").Value = -16710
$ws.Range("H138").Value = 40912.5
$ws.Range("J138").Value = 40912.5
$ws.Range("L138").Value = 40912.5
$ws.Range("N138").Value = -51192.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 45131.6
$ws.Range("J20").Value = 45131.6
$ws.Range("L20").Value = 45131.6
$ws.Range("N20").Value = -45603.6
$ws.Range("H30").Value = 45131.6
$ws.Range("J30").Value = 45131.6
$ws.Range("L30").Value = 45131.6
$ws.Range("N30").Value = -45313.6
$ws.Range("H31").Value = 3662.5908
$ws.Range("I31").Value = 1084.7778
$ws.Range("K31").Value = 1084.7778
$ws.Range("M31").Value = -789.7778000000001
$ws.Range("H34").Value = 3662.5908
$ws.Range("I34").Value = 1084.7778
$ws.Range("K34").Value = 1084.7778
$ws.Range("M34").Value = -882.7778000000001
$ws.Range("H99").Value = 16671115
$ws.Range("I99").Value = 25001798
$ws.Range("J99").Value = 9750
$ws.Range("K99").Value = 25001798
$ws.Range("L99").Value = 9750
$ws.Range("M99").Value = -25000300
$ws.Range("N99").Value = -12746
$ws.Range("H126").Value = 16671115
$ws.Range("I126").Value = 25001798
$ws.Range("J126").Value = 9750
$ws.Range("K126").Value = 75005394
$ws.Range("L126").Value = 29250
$ws.Range("M126").Value = -75002924
$ws.Range("N126").Value = -34190
$ws.Range("H127").Value = 41850
$ws.Range("J127").Value = 41850
$ws.Range("L127").Value = 41850
$ws.Range("N127").Value = -51770
$ws.Range("H128").Value = 45131.6
$ws.Range("J128").Value = 45131.6
$ws.Range("L128").Value = 45131.6
$ws.Range("N128").Value = -55091.6
$ws.Range("H130").Value = 41864
$ws.Range("J130").Value = 41864
$ws.Range("L130").Value = 41864
$ws.Range("N130").Value = -51904
$ws.Range("H138").Value = 43411.11
$ws.Range("J138").Value = 43411.11
$ws.Range("L138").Value = 43411.11
$ws.Range("N138").Value = -53691.11
$ws.Range("H140").Value = 88046.664
$ws.Range("J140").Value = 88046.664
$ws.Range("L140").Value = 88046.664
$ws.Range("N140").Value = -98406.664
$ws.Range("H141").Value = 21859.273
$ws.Range("J141").Value = 21859.273
$ws.Range("L141").Value = 21859.273
$ws.Range("N141").Value = -32219.273

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 98.75
$ws.Range("I12").Value = 13.333333
$ws.Range("K12").Value = 39.999999
$ws.Range("M12").Value = 133.000001
$ws.Range("H131").Value = 769.66
$ws.Range("I131").Value = 340
$ws.Range("J131").Value = 792.2737
$ws.Range("K131").Value = 1020
$ws.Range("L131").Value = 2376.8211
$ws.Range("M131").Value = 4020
$ws.Range("N131").Value = -12456.8211
$ws.Range("H140").Value = 4571.1113
$ws.Range("I140").Value = 4986.25
$ws.Range("J140").Value = 1250
$ws.Range("K140").Value = 14958.75
$ws.Range("L140").Value = 3750
$ws.Range("M140").Value = -9778.75
$ws.Range("N140").Value = -14110

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 8184045.5
$ws.Range("I10").Value = 10001167
$ws.Range("J10").Value = 6999.5
$ws.Range("K10").Value = 10001167
$ws.Range("L10").Value = 6999.5
$ws.Range("M10").Value = -10000998
$ws.Range("N10").Value = -7337.5
$ws.Range("H12").Value = 21517
$ws.Range("J12").Value = 29774.75
$ws.Range("L12").Value = 29774.75
$ws.Range("N12").Value = -30054.75
$ws.Range("H132").Value = 5789.8
$ws.Range("I132").Value = 3000
$ws.Range("K132").Value = 9000
$ws.Range("M132").Value = -6470
$ws.Range("H140").Value = 37961.43
$ws.Range("J140").Value = 37961.43
$ws.Range("L140").Value = 37961.43
$ws.Range("N140").Value = -48321.43

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5810.625
$ws.Range("I7").Value = 3667.6667
$ws.Range("J7").Value = 7096.4
$ws.Range("K7").Value = 3667.6667
$ws.Range("L7").Value = 7096.4
$ws.Range("M7").Value = -3555.6667
$ws.Range("N7").Value = -7320.4
$ws.Range("H40").Value = 9325.429
$ws.Range("I40").Value = 9800
$ws.Range("J40").Value = 9246.333000000001
$ws.Range("K40").Value = 9800
$ws.Range("L40").Value = 9246.333000000001
$ws.Range("M40").Value = -9664
$ws.Range("N40").Value = -9518.333000000001
$ws.Range("H93").Value = 4631318.5
$ws.Range("J93").Value = 1818.7273
$ws.Range("L93").Value = 1818.7273
$ws.Range("N93").Value = -4314.7273
$ws.Range("H122").Value = 4564.8945
$ws.Range("I122").Value = 2359.8
$ws.Range("J122").Value = 5352.4287
$ws.Range("K122").Value = 7079.400000000001
$ws.Range("L122").Value = 16057.2861
$ws.Range("M122").Value = -4629.400000000001
$ws.Range("N122").Value = -20957.2861
$ws.Range("H126").Value = 5810.625
$ws.Range("I126").Value = 3667.6667
$ws.Range("J126").Value = 7096.4
$ws.Range("K126").Value = 11003.0001
$ws.Range("L126").Value = 21289.2
$ws.Range("M126").Value = -8533.000100000001
$ws.Range("N126").Value = -26229.2
$ws.Range("H132").Value = 12182.576
$ws.Range("I132").Value = 14041.85
$ws.Range("J132").Value = 9322.154
$ws.Range("K132").Value = 42125.55
$ws.Range("L132").Value = 27966.462
$ws.Range("M132").Value = -39595.55
$ws.Range("N132").Value = -33026.462
$ws.Range("H139").Value = 46481.668
$ws.Range("J139").Value = 46481.668
$ws.Range("L139").Value = 46481.668
$ws.Range("N139").Value = -56761.668
$ws.Range("H141").Value = 35296.367
$ws.Range("J141").Value = 35296.367
$ws.Range("L141").Value = 35296.367
$ws.Range("N141").Value = -45656.367

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 40652.75
$ws.Range("I23").Value = 30866.666
$ws.Range("J23").Value = 70011
$ws.Range("K23").Value = 30866.666
$ws.Range("L23").Value = 70011
$ws.Range("M23").Value = -30637.666
$ws.Range("N23").Value = -70469
$ws.Range("H126").Value = 2351.3333
$ws.Range("I126").Value = 1332.5
$ws.Range("K126").Value = 3997.5
$ws.Range("M126").Value = -1527.5
$ws.Range("H138").Value = 39759.6
$ws.Range("J138").Value = 39759.6
$ws.Range("L138").Value = 39759.6
$ws.Range("N138").Value = -50039.6
$ws.Range("H140").Value = 32771.6
$ws.Range("J140").Value = 32771.6
$ws.Range("L140").Value = 32771.6
$ws.Range("N140").Value = -43131.6
$ws.Range("H141").Value = 25854
$ws.Range("J141").Value = 25854
$ws.Range("L141").Value = 25854
$ws.Range("N141").Value = -36214
